$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'26.769.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.14%  "
$ws.Range("D3").Value = "'1.693.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.63%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'219.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.56%  "
$ws.Range("D6").Value = "'0.5121"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -13.47%  "
$ws.Range("D7").Value = "'1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.2572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.75%  "
$ws.Range("D9").Value = "'21.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.72%  "
$ws.Range("D10").Value = "'0.06164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.77%  "
$ws.Range("D11").Value = "'0.07338"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").Value = "'1.695.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.45%  "
$ws.Range("D13").Value = "'4.458"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.80%  "
$ws.Range("D14").Value = "'0.5769"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.77%  "
$ws.Range("D15").Value = "'1.923.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.66%  "
$ws.Range("D16").Value = "'0.000008184"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -11.91%  "
$ws.Range("D17").Value = "'65.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -12.79%  "
$ws.Range("D18").Value = "'26.778.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.49%  "
$ws.Range("D19").Value = "'5.029"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.05%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "'10.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.69%  "
$ws.Range("D22").Value = "'186.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -11.00%  "
$ws.Range("D23").Value = "'6.239"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.98%  "
$ws.Range("D24").Value = "'1.008"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'142.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.81%  "
$ws.Range("D26").Value = "'7.465"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.53%  "
$ws.Range("D27").Value = "'0.1142"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.14%  "
$ws.Range("D28").Value = "'15.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.22%  "
$ws.Range("D29").Value = "'1.332"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.46%  "
$ws.Range("D30").Value = "'0.05843"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.61%  "
$ws.Range("D31").Value = "'1.343"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.07%  "
$ws.Range("D32").Value = "'3.459"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.58%  "
$ws.Range("D33").Value = "'3.422"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.80%  "
$ws.Range("D34").Value = "'1.640"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").Value = "'0.9921"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.87%  "
$ws.Range("E36").Value = "  -4.22%  "
$ws.Range("D37").Value = "'0.5947"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.36%  "
$ws.Range("D38").Value = "'2.667"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("D39").Value = "'1.087.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("D40").Value = "'0.01589"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.25%  "
$ws.Range("D41").Value = "'0.8558"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").Value = "'5.824"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.66%  "
$ws.Range("D43").Value = "'1.006"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").Value = "'97.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.92%  "
$ws.Range("E45").Value = "  -6.48%  "
$ws.Range("D46").Value = "'55.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.65%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'0.00000000104"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.48%  "
$ws.Range("D49").Value = "'8.030"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.4329"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05236"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.95%  "
